$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values like "1.00" or "0.0000254" would otherwise be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '66.212.85'
$ws.Range("E2").Value = '  -4.36%  '

# Row 3
$ws.Range("D3").Value = '3.558.91'
$ws.Range("E3").Value = '  -4.83%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").Value = '583.66'
$ws.Range("E5").Value = '  -5.19%  '

# Row 6
$ws.Range("D6").Value = '186.11'
$ws.Range("E6").Value = '  +0.44%  '

# Row 7
$ws.Range("D7").Value = '3.550.77'
$ws.Range("E7").Value = '  -5.05%  '

# Row 8
$ws.Range("D8").Value = '0.611'
$ws.Range("E8").Value = '  -4.57%  '

# Row 9
$ws.Range("E9").Value = '  +0.20%  '

# Row 10
$ws.Range("D10").Value = '0.665'
$ws.Range("E10").Value = '  -8.06%  '

# Row 11
$ws.Range("D11").Value = '0.145'
$ws.Range("E11").Value = '  -11.38%  '

# Row 12
$ws.Range("D12").Value = '53.33'
$ws.Range("E12").Value = '  -8.32%  '

# Row 13
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  -14.34%  '

# Row 14
$ws.Range("D14").Value = '9.68'
$ws.Range("E14").Value = '  -9.63%  '

# Row 15
$ws.Range("D15").Value = '4.122.88'
$ws.Range("E15").Value = '  -4.92%  '

# Row 16
$ws.Range("D16").Value = '3.546.77'
$ws.Range("E16").Value = '  -4.97%  '

# Row 17
$ws.Range("E17").Value = '  -0.92%  '

# Row 18
$ws.Range("D18").Value = '18.18'
$ws.Range("E18").Value = '  -7.12%  '

# Row 19
$ws.Range("D19").Value = '12.13'
$ws.Range("E19").Value = '  -6.93%  '

# Row 20
$ws.Range("D20").Value = '65.974.74'
$ws.Range("E20").Value = '  -4.42%  '

# Row 21
$ws.Range("D21").Value = '1.05'
$ws.Range("E21").Value = '  -7.94%  '

# Row 22
$ws.Range("D22").Value = '393.16'
$ws.Range("E22").Value = '  -5.49%  '

# Row 23
$ws.Range("D23").Value = '4.34'
$ws.Range("E23").Value = '  -7.01%  '

# Row 24
$ws.Range("D24").Value = '84.98'
$ws.Range("E24").Value = '  -5.42%  '

# Row 25
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '12.38'
$ws.Range("E25").Value = '  -3.42%  '

# Row 26
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '2.85'
$ws.Range("E26").Value = '  -6.96%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '10.68'
$ws.Range("E27").Value = '  -3.09%  '

# Row 28
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").Value = '6.03'
$ws.Range("E28").Value = '  -1.01%  '

# Row 29
$ws.Range("D29").Value = '3.52'
$ws.Range("E29").Value = '  -8.08%  '

# Row 30
$ws.Range("D30").Value = '8.88'
$ws.Range("E30").Value = '  -7.89%  '

# Row 31
$ws.Range("D31").Value = '30.70'
$ws.Range("E31").Value = '  -8.00%  '

# Row 32
$ws.Range("D32").Value = '6.79'
$ws.Range("E32").Value = '  -8.02%  '

# Row 33
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '617.56'
$ws.Range("E33").Value = '  -1.33%  '

# Row 34
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").Value = '12.06'
$ws.Range("E34").Value = '  -5.60%  '

# Row 35
$ws.Range("D35").Value = '63.06'
$ws.Range("E35").Value = '  -4.66%  '

# Row 36
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  -9.08%  '

# Row 37
$ws.Range("D37").Value = '41.05'
$ws.Range("E37").Value = '  -8.48%  '

# Row 38
$ws.Range("E38").Value = '  -0.01%  '

# Row 39
$ws.Range("D39").Value = '0.374'
$ws.Range("E39").Value = '  -7.64%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0740'
$ws.Range("E40").Value = '  -16.86%  '

# Row 41
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.20%  '

# Row 42
$ws.Range("D42").Value = '0.130'
$ws.Range("E42").Value = '  -9.17%  '

# Row 43
$ws.Range("D43").Value = '2.933.51'
$ws.Range("E43").Value = '  +3.76%  '

# Row 44
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  -10.53%  '

# Row 45
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '2.44'
$ws.Range("E45").Value = '  -8.05%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0403'
$ws.Range("E46").Value = '  -9.30%  '

# Row 47
$ws.Range("D47").Value = '3.08'
$ws.Range("E47").Value = '  -3.61%  '

# Row 48
$ws.Range("E48").Value = '  -6.31%  '

# Row 49
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '8.53'
$ws.Range("E49").Value = '  -8.05%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '138.20'
$ws.Range("E50").Value = '  -2.49%  '

# Row 51
$ws.Range("D51").Value = '2.76'
$ws.Range("E51").Value = '  -0.38%  '
